# Applies:
#  - Widen a few columns on the "semantic_aspect_model_schema" sheet
#  - Replace single-underscore flattened-field separators with a
#    double-underscore separator ("type_certificateType" -> "type__certificateType", etc.)
#    both in the schema header row and in the "description" sheet's legend rows.
#  - Renumber the "digital twin fields" legend bullet from "2." to "1."

$wb = $excel.ActiveWorkbook

$schemaSheet = $wb.Worksheets.Item("semantic_aspect_model_schema")
$descriptionSheet = $wb.Worksheets.Item("description")

# --- Column width updates on the schema sheet ---
# NOTE: the host snaps ColumnWidth to the nearest 1/6-character pixel grid
# (it internally adds ~5px of padding and rounds), so the raw target widths
# below are pre-compensated (target - 5/6, rounded to the 1/6 grid) so the
# stored <col width="…"> comes out as close as possible to the intended
# 25.2 / 28.8 / 39.6 / 42 / 28.8 / 27.6 values.
$schemaSheet.Columns.Item(2).ColumnWidth = 24.333333333333332   # B: 24 -> 25.2
$schemaSheet.Columns.Item(3).ColumnWidth = 28                   # C: 27.6 -> 28.8
$schemaSheet.Columns.Item(6).ColumnWidth = 38.833333333333336   # F: 38.4 -> 39.6
$schemaSheet.Columns.Item(7).ColumnWidth = 41.166666666666664   # G: 40.8 -> 42
$schemaSheet.Columns.Item(12).ColumnWidth = 28                  # L: 27.6 -> 28.8
$schemaSheet.Columns.Item(13).ColumnWidth = 26.833333333333332  # M: 26.4 -> 27.6

# --- Header row field-name separator fix on the schema sheet ---
$schemaSheet.Range("B1").Value = "type__certificateType"
$schemaSheet.Range("C1").Value = "type__certificateVersion"
$schemaSheet.Range("F1").Value = "enclosedSites[0]__enclosedSiteBpn"
$schemaSheet.Range("G1").Value = "enclosedSites[0]__areaOfApplication"
$schemaSheet.Range("L1").Value = "validator__validatorName"
$schemaSheet.Range("M1").Value = "validator__validatorBpn"

# --- "description" sheet updates ---
$descriptionSheet.Range("A3").Value = "1. Columns highlighted in olive green are digital twin fields."

$descriptionSheet.Range("A6").Value = "type__certificateType"
$descriptionSheet.Range("A7").Value = "type__certificateVersion"
$descriptionSheet.Range("A10").Value = "enclosedSites[0]__enclosedSiteBpn"
$descriptionSheet.Range("A11").Value = "enclosedSites[0]__areaOfApplication"
$descriptionSheet.Range("A16").Value = "validator__validatorName"
$descriptionSheet.Range("A17").Value = "validator__validatorBpn"
